# "E suite.xlsx" / "Test Cases" sheet
#
# The diff flips the Runmode for the two later watch-list test cases
# (TestCase_E2 / row 3 and TestCase_E3 / row 4) from "Y" to "N", and moves
# the view's scroll/selection to column B with C8 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 3 = TestCase_E2 ("...watchlist from document page"), Runmode column C.
$ws.Range("C3").Value = "N"

# Row 4 = TestCase_E3 ("...delete a document from watchlist"), Runmode column C.
$ws.Range("C4").Value = "N"

# Make "Test Cases" the active sheet and restore the view's selection state
# (scrolled so column B is left-most, active cell C8).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C8").Select()
